$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy the formatting from the neighboring
# header cell (G1) so it matches the other header cells' style, then set
# the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data value for the new column in row 2
$ws.Range("H2").Value = 0
